$d = $word.ActiveDocument

# Replace the three-digit / one-digit division answers throughout the document.
# Each (old -> new) pair below is unique in the document, so a simple
# Find/Replace over the whole content is safe and order-independent.
$replacements = @(
    @('711÷6=118, 3', '422÷7=60, 2'),
    @('272÷7=38, 6', '100÷6=16, 4'),
    @('528÷5=105, 3', '112÷2=56, 0'),
    @('129÷6=21, 3', '311÷5=62, 1'),
    @('180÷9=20, 0', '719÷7=102, 5'),
    @('317÷6=52, 5', '838÷8=104, 6'),
    @('532÷3=177, 1', '373÷4=93, 1'),
    @('308÷4=77, 0', '853÷4=213, 1'),
    @('729÷2=364, 1', '606÷2=303, 0'),
    @('895÷8=111, 7', '266÷3=88, 2'),
    @('237÷6=39, 3', '775÷4=193, 3'),
    @('553÷2=276, 1', '128÷6=21, 2'),
    @('313÷4=78, 1', '384÷2=192, 0'),
    @('419÷6=69, 5', '753÷8=94, 1'),
    @('479÷5=95, 4', '764÷5=152, 4'),
    @('371÷5=74, 1', '979÷2=489, 1'),
    @('247÷6=41, 1', '320÷7=45, 5'),
    @('690÷3=230, 0', '754÷3=251, 1'),
    @('962÷8=120, 2', '272÷4=68, 0'),
    @('178÷2=89, 0', '547÷7=78, 1'),
    @('330÷3=110, 0', '796÷5=159, 1'),
    @('510÷4=127, 2', '458÷3=152, 2'),
    @('410÷3=136, 2', '520÷9=57, 7'),
    @('420÷2=210, 0', '875÷3=291, 2'),
    @('720÷9=80, 0', '583÷6=97, 1'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Replacements applied: $($replacements.Count)"
